$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the shared video-path string (index 16) by setting all three
#    existing references to the same new value so Excel mutates the shared
#    string in place rather than allocating a new one.
$videoPath = '/videos/saiareact.mp4'
$ws.Cells.Item(3,3).Value = $videoPath
$ws.Cells.Item(4,3).Value = $videoPath
$ws.Cells.Item(6,3).Value = $videoPath

# 2) Adjust row heights for the existing program rows (row 4 keeps its
#    original 409.5 max height).
$ws.Rows.Item(3).RowHeight = 360
$ws.Rows.Item(5).RowHeight = 405
$ws.Rows.Item(6).RowHeight = 360

# 3) Column B needs to be widened to fit the new, longer program titles.
$ws.Columns.Item(2).ColumnWidth = 36.1667

# 4) Append the 21 new program rows (7-27). Formats are cloned from
#    existing rows so the new cells pick up the same wrap-text / font
#    styles used throughout the table, then the text values + row
#    heights are applied on top.

# Row 7
$ws.Range("A6:D6").Copy()
$ws.Range("A7:D7").PasteSpecial(-4122)
$ws.Cells.Item(7,1).Value = 'Culinary Arts'
$ws.Cells.Item(7,2).Value = 'Chop, fillet and create exciting cuisine every day as you experience the fast paced culinary industry.  From cooking for the student run Spartan Inn to providing food for large groups, students gain valuable restaurant experience.  Bake breads, cakes, pies, experiment with confections and chocolate.'
$ws.Cells.Item(7,4).Value = '/Culinary'
$ws.Rows.Item(7).RowHeight = 135

# Row 8
$ws.Range("A6:D6").Copy()
$ws.Range("A8:D8").PasteSpecial(-4122)
$ws.Cells.Item(8,1).Value = 'Architecture & Construction Technology'
$ws.Cells.Item(8,2).Value = '
The Architecture and Construction Technology program prepares individuals to apply technical knowledge and skills to the basic aspects of planning, preparing and interpreting architectural, structural, electrical, topographical and other drawings and sketches used in various Architectural and Construction fields.
Instruction is designed to provide experiences in AutoCAD and Revit design software for production of Architectural drawings. Instruction will be provided in pro-estimating and fastrack software for estimating and scheduling.
In addition to digital modeling we will also be creating physical models to understand special relationships. Students will have the opportunity to be dual enrolled their Junior and Senior year where they will obtain 18 credits with H.A.C.C. towards their Architecture or Construction Management degrees.'
$ws.Cells.Item(8,4).Value = '/Architecture'
$ws.Rows.Item(8).RowHeight = 409.5

# Row 9
$ws.Range("A6:D6").Copy()
$ws.Range("A9:D9").PasteSpecial(-4122)
$ws.Cells.Item(9,1).Value = 'Automotive Technology '
$ws.Cells.Item(9,2).Value = 'Automotive Technology offers a wide variety of opportunities for the students.  They can learn a valuable trade through a combination of classroom instruction and hands-on training.  Curriculum is revised to keep current with the latest technological developments in the automotive industry. 
Some key areas of instruction are:  Automotive maintenance and light repair, engine repair, automatic transmission/transaxle, manual drive train and axle, suspension and steering, brakes, electrical/electronic systems, heating and air conditioning, and engine performance.'
$ws.Cells.Item(9,4).Value = '/Automotive'
$ws.Rows.Item(9).RowHeight = 270

# Row 10
$ws.Range("A6:D6").Copy()
$ws.Range("A10:D10").PasteSpecial(-4122)
$ws.Cells.Item(10,1).Value = 'Carpentry & Construction Management'
$ws.Cells.Item(10,2).Value = 'The Carpentry program covers the major aspects of residential and light commercial construction as well as a variety of remodeling functions. Students will work with a variety of building and finishing materials in order to become familiar with proper methods. Training is given using hand tools, power tools, and other portable equipment. Onsite and offsite construction projects are incorporated into classroom instruction.'
$ws.Cells.Item(10,3).Value = ' '
$ws.Cells.Item(10,4).Value = '/Carpentry'
$ws.Rows.Item(10).RowHeight = 180

# Row 11
$ws.Range("A6:D6").Copy()
$ws.Range("A11:D11").PasteSpecial(-4122)
$ws.Cells.Item(11,1).Value = 'Collision Repair Technology'
$ws.Cells.Item(11,2).Value = 'This is a high tech area that is designed to provide training for persons interested to provide training for persons interested in learning how to repair damaged vehicles.  Students will learn how to repair damaged vehicles.  Students will learn how to repair wrecked vehicles using the latest techniques and equipment available.  This technical area features a Chief E-Z Liner frame machine, Dupont paint system, the latest spray guns available, and an Omni down draft spray booth with an oven.  Hands-on activities are mixed with the correct blend of theory to allow the students to master the techniques of frame repair, MIG Welding, painting and custom designs, dent removal, suspension systems, brakes, electrical systems, and state inspection.'
$ws.Cells.Item(11,4).Value = '/Collision'
$ws.Rows.Item(11).RowHeight = 315

# Row 12
$ws.Range("A6:D6").Copy()
$ws.Range("A12:D12").PasteSpecial(-4122)
$ws.Cells.Item(12,1).Value = 'Commercial & Advertising Art'
$ws.Cells.Item(12,2).Value = 'Students will prepare a portfolio throughout the program to promote their work and talent when they graduate.  The major emphasis is on the basic principles of design and elements of art through skill development and exploring different media.  Special emphasis is placed on manual as well as computer illustration, layout, composition, and photography skills.  Students will prepare graphic and advertising projects from the idea stage through to web/pre-press.  The students will be using industry software throughout this course learning the basic skills to advance.'
$ws.Cells.Item(12,4).Value = '/Commerical'
$ws.Rows.Item(12).RowHeight = 240

# Row 13
$ws.Range("A6:D6").Copy()
$ws.Range("A13:D13").PasteSpecial(-4122)
$ws.Cells.Item(13,1).Value = 'Communications Technology'
$ws.Cells.Item(13,2).Value = 'Communication Technology is a multimedia course preparing students for careers and post-secondary education in digital and print fields.  Explore audio and video media using state of the art studio equipment, cameras, lighting, and live production visual imaging software.  Create, print, and assemble publications, design and print t-shirts, and perform tasks related to desktop publishing.  Students will learn and apply industry relevant computer software applications to real world media projects in a wide variety of mediums.'
$ws.Cells.Item(13,4).Value = '/Communications'
$ws.Rows.Item(13).RowHeight = 225

# Row 14
$ws.Range("A6:D6").Copy()
$ws.Range("A14:D14").PasteSpecial(-4122)
$ws.Cells.Item(14,1).Value = 'Criminal Justice'
$ws.Cells.Item(14,2).Value = 'Do you have what it takes?  Personal responsibility, integrity, and self discipline are needed to be successful in the field of public safety.  Improve your physical training, take fingerprints, and perform first aid.  Make quick decisions while using the Fire Arms Training Simulator, and the Geospacial Information Systems mapping software. This program prepares students for entry level positions in the Criminal Justice arena, specifically, the field of Law Enforcement and Emergency Medical Services.  Personal responsibility, integrity, and self-discipline will be expected.  Intensive physical training, rigorous academics, and adherence to strict standards of personal appearance are required.'
$ws.Cells.Item(14,4).Value = '/Criminal '
$ws.Rows.Item(14).RowHeight = 285

# Row 15
$ws.Range("A6:B6").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Cells.Item(15,1).Value = 'Dental Assisting'
$ws.Cells.Item(15,2).Value = 'Work with patients in a dental office and work with the dental team. Learn to help the dentist chairside with procedures by passing dental instruments and materials to the dentist and keeping the mouth dry.   Mix and prepare dental materials, sterilize instruments, disinfect treatment rooms, take impressions and take x-rays.  Perform basic lab duties such as and make models of teeth. '
$ws.Cells.Item(15,4).Value = '/Dental'
$ws.Rows.Item(15).RowHeight = 165

# Row 16
$ws.Range("A6:B6").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Cells.Item(16,1).Value = 'Diesel Technology'
$ws.Cells.Item(16,2).Value = ' Diesel Technology is a broad-based program covering basic diesel engine fundamentals, full truck alignment, electrical/electronic engines, precision measuring, engine overhaul, air brake systems, drive train components, and PA State inspection procedures.  This program also offers state of the art equipment which will be an asset in furthering your education. '
$ws.Cells.Item(16,4).Value = '/Diesel'
$ws.Rows.Item(16).RowHeight = 150

# Row 17
$ws.Range("A6:B6").Copy()
$ws.Range("A17:B17").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D17").PasteSpecial(-4122)
$ws.Cells.Item(17,1).Value = 'Early Childhood Education'
$ws.Cells.Item(17,2).Value = 'Are you creative, energetic, and enjoy working with children?  Learn how to develop and present lesson plans to young children, spend time working at local elementary schools, at our on-site Early Learning Center, or at other community childcare programs.  Earn hours towards the 480 hours of training required to apply for your Child Development Associate Credential.'
$ws.Cells.Item(17,4).Value = '/Childhood'
$ws.Rows.Item(17).RowHeight = 150

# Row 18
$ws.Range("A6:B6").Copy()
$ws.Range("A18:B18").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Cells.Item(18,1).Value = 'Electrical Occupations'
$ws.Cells.Item(18,2).Value = '
Safely install electrical systems in residential and commercial industrial buildings.  Read and understand blueprints.  Wire a variety of receptacles and outlet boxes.  Mount electrical panels.  Study the theory of electricity and learn how to operate and repair electric motors and electro-mechanical controls.'
$ws.Cells.Item(18,4).Value = '/Electrical'
$ws.Rows.Item(18).RowHeight = 150

# Row 19
$ws.Range("A6:B6").Copy()
$ws.Range("A19:B19").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D19").PasteSpecial(-4122)
$ws.Cells.Item(19,1).Value = 'Engineering/Advanced Manufacturing'
$ws.Cells.Item(19,2).Value = 'This technical program prepares students to apply knowledge and skills in the engineering field.  Basic instruction is provided in a variety of areas associated with engineering such as civil engineering, electrical and electronic engineering, electromechanical instrumentation, industrial production and mechanical engineering.  Instruction includes but is not limited to electrical circuitry, electronic digital and microprocessor applications, high and low voltage applications, instrumentation calibration, prototype development, testing, inspecting, systems analysis and maintenance, applications to specific engineering systems, CAD/CAM, fluid power, heating and cooling, manufacturing systems, principles of mechanics, properties of materials, and report writing.'
$ws.Cells.Item(19,4).Value = '/Engineering'
$ws.Rows.Item(19).RowHeight = 330

# Row 20
$ws.Range("A6:B6").Copy()
$ws.Range("A20:B20").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Cells.Item(20,1).Value = 'Heating, Ventilation & Air Conditioning/Plumbing'
$ws.Cells.Item(20,2).Value = '
The HVAC/Plumbing program trains students for the always changing, always in demand HVAC industry.  Students will learn to identify and safety use power and hand tools to complete hands-on projects.  They advanced to learning the materials and equipment of the field while completing more challenging installation jobs.  Students receive instruction on electricity and electrical components and controls, and the technical aspects of HVAC which includes proper installation, setup, service, and troubleshooting of residential and commercial heating and cooling systems.'
$ws.Cells.Item(20,4).Value = '/HVAC'
$ws.Rows.Item(20).RowHeight = 270

# Row 21
$ws.Range("A6:B6").Copy()
$ws.Range("A21:B21").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D21").PasteSpecial(-4122)
$ws.Cells.Item(21,1).Value = 'Landscaping & Agriculture'
$ws.Cells.Item(21,2).Value = 'The Landscaping & Agriculture Program covers training in equipment usage such as skid loader and Bobcat utility vehicle, greenhouse management, turf maintenance, plan identification, tree maintenance, nursery production and landscape design, management, and construction including pond and patio installation.  Practical experience is gained through school landscape plantings, plan and maintenance of a vegetable garden and maintenance of a nursery area.  Design is also taught using the Punch! Landscape Deck and Patio Designer program.'
$ws.Cells.Item(21,4).Value = '/Landscaping'
$ws.Rows.Item(21).RowHeight = 225

# Row 22
$ws.Range("A6:B6").Copy()
$ws.Range("A22:B22").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Cells.Item(22,1).Value = 'Licensed Cosmetology'
$ws.Cells.Item(22,2).Value = 'Students in the Licensed Cosmetology program study topics such as:  safety, history, careers, professionalism, infection control, chemistry, PA State Cosmetology law, salon business, and anatomy & physiology.  In addition, students are training in hair care skills such as shampooing, hair cutting, styling, braiding, perming, relaxing, coloring, foiling, extensions, and up-dos.  Instruction also includes the care of hands and nails, feet and toes, skin and makeup artistry.  Students practice these techniques on mannequins, practice hands, or people. 
During their junior and senior year, students participate in clinical experience by demonstrating services on customers at the cosmetology clinic-“Spartan Spa”. 
Upon completion of 1250 hours, a transcript of 75% of higher in the course, and an Official Criminal Record History check, students are able to take the State Board Examination, which is required for licensing.
'
$ws.Cells.Item(22,4).Value = '/Cosmetology'
$ws.Rows.Item(22).RowHeight = 409.5

# Row 23
$ws.Range("A6:B6").Copy()
$ws.Range("A23:B23").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Cells.Item(23,1).Value = 'Medical Professions Program'
$ws.Cells.Item(23,2).Value = '
Students are required to complete the foundation courses during the first two years at York Tech.  The intent of this program is to provide a foundation for these students and to aid students and teachers in choosing appropriate pathway placements that will lead to informed career decisions.
The Medical Professions Program includes several foundation courses to help students succeed in this high-paced, challenging career.  During the first two years, students will cover the following topics:
    Orientation & Safety
    Nutrition & Hydration
    Medical Terminology/Human Body
    Rehabilitation & Restoration
    Diseases and Disorders
    Basic Clinical Skills
    Mathematics in Allied Health
    Legal & Ethical Issues
    Emergency Care & Disaster Preparedness'
$ws.Cells.Item(23,4).Value = '/Medical'
$ws.Rows.Item(23).RowHeight = 405

# Row 24
$ws.Range("A6:B6").Copy()
$ws.Range("A24:B24").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D24").PasteSpecial(-4122)
$ws.Cells.Item(24,1).Value = 'PowerSports & Engine Technology'
$ws.Cells.Item(24,2).Value = 'PowerSports and Engine Technology is designed to provide training for persons interested in repairing PowerSports equipment.  The areas that will be covered include:  motorcycles, marine products, lift trucks, lawn maintenance products, and the engines that power them.  '
$ws.Cells.Item(24,4).Value = '/Powersports'
$ws.Rows.Item(24).RowHeight = 120

# Row 25
$ws.Range("A6:B6").Copy()
$ws.Range("A25:B25").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Cells.Item(25,1).Value = 'Precision Machining Technology'
$ws.Cells.Item(25,2).Value = 'This program builds a firm foundation in the following areas:  lathe, mill, grinder, drill press, bench work, precision measurement, print reading, and understanding geometric tolerances.  Students cover advanced machining techniques on Computer Numerical Control machines (CNC).'
$ws.Cells.Item(25,4).Value = '/PrecisionMachining'
$ws.Rows.Item(25).RowHeight = 120

# Row 26
$ws.Range("A6:B6").Copy()
$ws.Range("A26:B26").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Cells.Item(26,1).Value = 'Sports Technology & Exercise Science'
$ws.Cells.Item(26,2).Value = '
Do you enjoy physical activity and helping others?  Can you see yourself being a personal trainer, athletic trainer or physical therapist?  Learn the basic skills associated with athletic health care, basic anatomy and physiology.  Tape, wrap and brace athletic injuries.  Understand the proper ways to use heat and ice.  Chart vital signs such as blood pressure, pulse, temperature, and pain.'
$ws.Cells.Item(26,4).Value = '/Sports'
$ws.Rows.Item(26).RowHeight = 180

# Row 27
$ws.Range("A6:B6").Copy()
$ws.Range("A27:B27").PasteSpecial(-4122)
$ws.Range("D6").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Cells.Item(27,1).Value = 'Welding & Metal Fabrication Technology'
$ws.Cells.Item(27,2).Value = '
Welding/Metal Fabrication Technology program is a very diverse area in which the students are trained on Shielded Metal Arc Welding, Flux Cored Arc Welding, Gas Tungsten Arc welding and Oxy-Acetylene/Cutting, Welding and Brazing processes. Students also study emerging technologies such as glass and plastic welding.  Theory includes the aspects relating to safety in today''s welding facilities as well as welding symbol and blueprint usage.'
$ws.Cells.Item(27,4).Value = '/Welding'
$ws.Rows.Item(27).RowHeight = 195

$excel.CutCopyMode = $false

# 5) Update the view: scroll/selection now centers on row 11 as in the target file.
$win = $excel.ActiveWindow
$win.ScrollRow = 11
$win.ScrollColumn = 1
$ws.Range("G11").Select()
